$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Coin name) changes - rows whose rank order changed
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("B25").Value = "Monero"
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("B44").Value = "TrustWalletToken"

# Column C (Link) changes - rows whose rank order changed
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"

# Column D (Price) changes - force Text format so values like "214.24" are not
# auto-converted to floating point numbers by Excel
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.838.03"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.35"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.24"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0632"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.91"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.630.29"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.52"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.844.09"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.38"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.96"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.01"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.136.12"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.46"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0156"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.763.92"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.13"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"

# Column E (Volume 1h) changes
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +0.79%  "
